$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Ayo Dosunmu", "PG,SG,SF", "Chicago Bulls"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Jalen Suggs", "PG,SG", "Orlando Magic"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Malcolm Brogdon", "PG,SG", "Washington Wizards"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Dejounte Murray", "PG,SG", "New Orleans Pelicans"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jakob Poeltl", "C", "Toronto Raptors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
